# Update countries & provincias Spain
# Applies the 24-Sep-2020 13:48 data refresh to the "Pais" sheet:
#  - Nepal overtakes Costa Rica in the ranking (rows 54/55 swap country + values)
#  - San Pedro y Miquelon overtakes Groenlandia and Islas Malvinas (rows 213/214/216)
#  - Refreshed case counts for several other countries (no rank change)
#  - Updated "last refreshed" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 13:48"

# --- Row 16: Iran (rank unchanged) ---
$ws.Range("B16").Value = 436319
$ws.Range("C16").Value = 3521
$ws.Range("D16").Value = 367829
$ws.Range("E16").Value = 43475
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 175
$ws.Range("H16").Value = 25015

# --- Row 40: Kuwait (rank unchanged) ---
$ws.Range("B40").Value = 101851
$ws.Range("C40").Value = 552
$ws.Range("D40").Value = 92961
$ws.Range("E40").Value = 8298
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 592

# --- Row 49: Bielorrusia (rank unchanged) ---
$ws.Range("B49").Value = 76651
$ws.Range("C49").Value = 294
$ws.Range("D49").Value = 73733
$ws.Range("E49").Value = 2116
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 802

# --- Rows 54/55: Nepal overtakes Costa Rica ---
# Row 54 becomes Nepal with refreshed data
$ws.Range("A54").Value = "Nepal"
$ws.Range("B54").Value = 69301
$ws.Range("C54").Value = 1497
$ws.Range("D54").Value = 50411
$ws.Range("E54").Value = 18437
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 17
$ws.Range("H54").Value = 453

# Row 55 becomes Costa Rica, taking the old (pre-refresh) Nepal row's slot with
# the figures Costa Rica had before the refresh
$ws.Range("A55").Value = "Costa Rica"
$ws.Range("B55").Value = 68059
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 26136
$ws.Range("E55").Value = 41142
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 781

# --- Row 61: Suiza (rank unchanged) ---
$ws.Range("B61").Value = 51492
$ws.Range("C61").Value = 391
$ws.Range("D61").Value = 42300
$ws.Range("E61").Value = 7131
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 2061

# --- Row 88: Madagascar (rank unchanged) ---
$ws.Range("B88").Value = 16191
$ws.Range("C88").Value = 24
$ws.Range("D88").Value = 14833
$ws.Range("E88").Value = 1131
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 227

# --- Row 122: Hong Kong (rank unchanged) ---
$ws.Range("D122").Value = 4758
$ws.Range("E122").Value = 195

# --- Row 145: Malta (rank unchanged) ---
$ws.Range("B145").Value = 2898
$ws.Range("C145").Value = 42
$ws.Range("D145").Value = 2191
$ws.Range("E145").Value = 680
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 27

# --- Row 149: Islandia (rank unchanged) ---
$ws.Range("B149").Value = 2512
$ws.Range("C149").Value = 36
$ws.Range("D149").Value = 2150
$ws.Range("E149").Value = 352

# --- Rows 213/214/216: San Pedro y Miquelon overtakes Groenlandia and Islas Malvinas ---
# Row 213 becomes San Pedro y Miquelon with refreshed data
$ws.Range("A213").Value = "San Pedro y Miquelon"
$ws.Range("B213").Value = 16
$ws.Range("C213").Value = 4
$ws.Range("D213").Value = 6
$ws.Range("E213").Value = 10
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214 becomes Groenlandia, keeping its pre-refresh figures
$ws.Range("A214").Value = "Groenlandia"
$ws.Range("B214").Value = 14
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 14
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215 (Montserrat) is unchanged

# Row 216 becomes Islas Malvinas, keeping its pre-refresh figures
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

# Row 217 (Santa Sede) is unchanged
